# Simulated Wild Card round and logged it
# This script appends the results of one additional (road/"Road") playoff
# game - the Wild Card round - into the season-long tracking workbook:
#  - YDS sheet: appends per-play rush/pass yardage logs (OFF & DEF)
#  - OFF / DEF sheets: updates the Road row (and season totals) aggregate counts
#  - ST sheet: appends kickoff/punt distance & return logs, updates totals
#  - TURNS sheet: updates turnover totals
#  - PEN sheet: updates penalty totals

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly logged play yardages to the running lists
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 2 2 0 4 22 -5 2 3 3 4 2 3 18 8 -1 8 10 46 1 22 27 3 6 0 4 2 6 1 21 0 11 4 4 3 1"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 3 19 13 11 10 6 0 28 7 1 5 1 6 11 8 12"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 1 1 5 9 -1 2 2 -6 1 1 2 0 7 14 1 -1 3 17 -1 0 2 4 3 0 2 0 15"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 4 11 10 7 0 7 25 2 5 16 8 2 10 10 7 8 -2 8 20 11 3 6 10 6 3 11 14 1 11"

# ---------------------------------------------------------------------
# OFF sheet: Road-row (row 3) play counts plus season totals (row 2)
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 227
$wsOFF.Range("F2").Value = 64
$wsOFF.Range("G2").Value = 77
$wsOFF.Range("I2").Value = 8
$wsOFF.Range("J2").Value = 50
$wsOFF.Range("L2").Value = 314
$wsOFF.Range("M2").Value = 206
$wsOFF.Range("O2").Value = 27
$wsOFF.Range("P2").Value = 18
$wsOFF.Range("Q2").Value = 654

$wsOFF.Range("B3").Value = 11
$wsOFF.Range("C3").Value = 226
$wsOFF.Range("D3").Value = 4
$wsOFF.Range("E3").Value = 33
$wsOFF.Range("F3").Value = 138
$wsOFF.Range("H3").Value = 33
$wsOFF.Range("I3").Value = 68
$wsOFF.Range("J3").Value = 41
$wsOFF.Range("N3").Value = 27

# ---------------------------------------------------------------------
# DEF sheet: Road-row (row 3) play counts plus season totals (row 2)
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 188
$wsDEF.Range("F2").Value = 49
$wsDEF.Range("G2").Value = 52
$wsDEF.Range("J2").Value = 23
$wsDEF.Range("L2").Value = 318
$wsDEF.Range("M2").Value = 208
$wsDEF.Range("O2").Value = 18
$wsDEF.Range("P2").Value = 11
$wsDEF.Range("Q2").Value = 528

$wsDEF.Range("C3").Value = 200
$wsDEF.Range("E3").Value = 40
$wsDEF.Range("F3").Value = 111
$wsDEF.Range("G3").Value = 35
$wsDEF.Range("H3").Value = 34
$wsDEF.Range("I3").Value = 74
$wsDEF.Range("J3").Value = 48
$wsDEF.Range("N3").Value = 13

# ---------------------------------------------------------------------
# ST sheet: kickoff/punt totals, plus appended per-kick/punt logs
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 87
$wsST.Range("D2").Value = 69
$wsST.Range("F2").Value = 386
$wsST.Range("G2").Value = 382
$wsST.Range("J2").Value = 181
$wsST.Range("K2").Value = 178
$wsST.Range("L2").Value = 111
$wsST.Range("M2").Value = 100

$wsST.Range("B3").Value = 47

$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " 59 52 58"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " 32 13 21"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 19 23 27 26 18"
$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 44 40 36 23 50 39 36"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 7 0 0 0 23 0 0"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 0 0 0 0 0 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet: turnover totals
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B2").Value = 14
$wsTURNS.Range("C2").Value = 6
$wsTURNS.Range("D2").Value = 6
$wsTURNS.Range("E3").Value = 5

# ---------------------------------------------------------------------
# PEN sheet: penalty totals
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 18
$wsPEN.Range("B3").Value = 28
